# "Update for Magozzi resub"
# Rename the "SA" (South Africa?) site-code header/label strings to "EC"
# throughout the "ham" worksheet: both the header row (row 1, columns
# D:N) and the corresponding row labels (column A, rows 4:14) that use
# the same text.
#
#   OldSA.1_H_1 -> OldEC.1_H_1
#   OldSA.2_H_1 -> OldEC.2_H_1
#   SA_H_1..SA_H_9 -> EC_H_1..EC_H_9
#
# All other header/label text (DEN, OldUT, UT, CAN, US, VSMOW_H) is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNames = @("OldEC.1_H_1","OldEC.2_H_1","EC_H_1","EC_H_2","EC_H_3","EC_H_4","EC_H_5","EC_H_6","EC_H_7","EC_H_8","EC_H_9")

# Header row (row 1): columns D..N hold these 11 labels in order.
$headerCols = @("D","E","F","G","H","I","J","K","L","M","N")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $newNames[$i]
}

# Row labels (column A): rows 4..14 hold the same 11 labels in order.
for ($i = 0; $i -lt $newNames.Length; $i++) {
    $row = 4 + $i
    $ws.Range("A" + $row).Value = $newNames[$i]
}

# Move the active selection to A6 (was the whole sheet A1:XFD1048576).
[void]$ws.Range("A6").Select()
